{"js": "// Actualizaci\u00f3n del plan de pruebas\n// Rebuild the body with the new \"Sistema de Control de Asistencia\" test\n// plan content, dropping the old Heading1/Heading2 paragraph styles (the\n// new paragraphs are plain body text) and splitting what used to be\n// single runs containing <w:br/> line breaks into separate paragraphs,\n// with blank paragraphs between test cases.\n\nconst body = context.document.body;\n\n// Remember how many paragraphs exist today so we know what to remove\n// once the replacement content has been appended.\nconst oldParagraphs = body.paragraphs;\noldParagraphs.load(\"items\");\nawait context.sync();\nconst oldCount = oldParagraphs.items.length;\n\n// The new plan, one entry per paragraph (\"\" => a blank separator\n// paragraph). Appending via body.insertParagraph(..., \"End\") creates\n// plain paragraphs with no inherited paragraph style.\nconst newContent = [\n  \"PLAN DE PRUEBAS \u2013 Sistema de Control de Asistencia\",\n  \"\",\n  \"CP01: Registrar entrada v\u00e1lida\",\n  \"Entrada: Juan P\u00e9rez\",\n  \"Resultado esperado: Entrada registrada correctamente\",\n  \"\",\n  \"CP02: Registrar salida v\u00e1lida\",\n  \"Resultado esperado: Salida registrada correctamente\",\n  \"\",\n  \"CP03: Doble entrada\",\n  \"Resultado esperado: Error: ya existe una entrada sin salida\",\n  \"\",\n];\n\nfor (const line of newContent) {\n  body.insertParagraph(line, \"End\");\n}\nawait context.sync();\n\n// Remove the original paragraphs now that the new ones follow them.\nconst allParagraphs = body.paragraphs;\nallParagraphs.load(\"items\");\nawait context.sync();\nfor (let i = 0; i < oldCount; i++) {\n  allParagraphs.items[i].delete();\n}\nawait context.sync();\n", "ps1": "# Actualizaci\u00f3n del plan de pruebas\n# Rebuild the body with the new \"Sistema de Control de Asistencia\" test\n# plan content, dropping the old Heading1/Heading2 paragraph styles and\n# splitting the old <w:br/>-separated runs into their own paragraphs,\n# with blank paragraphs between test cases.\n\n$d = $word.ActiveDocument\n\n# Remember how many paragraphs exist today so we know what to remove\n# once the replacement content has been appended.\n$oldCount = $d.Paragraphs.Count\n\n# The new plan, one entry per paragraph (\"\" => a blank separator\n# paragraph). Paragraphs.Add() appends a plain paragraph at the end of\n# the document with no inherited paragraph style.\n$newContent = @(\n  \"PLAN DE PRUEBAS \u2013 Sistema de Control de Asistencia\",\n  \"\",\n  \"CP01: Registrar entrada v\u00e1lida\",\n  \"Entrada: Juan P\u00e9rez\",\n  \"Resultado esperado: Entrada registrada correctamente\",\n  \"\",\n  \"CP02: Registrar salida v\u00e1lida\",\n  \"Resultado esperado: Salida registrada correctamente\",\n  \"\",\n  \"CP03: Doble entrada\",\n  \"Resultado esperado: Error: ya existe una entrada sin salida\",\n  \"\"\n)\n\nforeach ($line in $newContent) {\n    $p = $d.Paragraphs.Add()\n    $p.Range.Text = $line\n}\n\n# Remove the original paragraphs now that the new ones follow them.\nfor ($i = $oldCount; $i -ge 1; $i--) {\n    $d.Paragraphs($i).Range.Delete()\n}\n"}
